$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.596.33"
$ws.Range("E2").Value = "'  +0.61%  "
$ws.Range("D3").Value = "'1.629.55"
$ws.Range("E3").Value = "'  +0.23%  "
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("D5").Value = "'213.17"
$ws.Range("E5").Value = "'  +0.22%  "
$ws.Range("E6").Value = "'  +2.63%  "
$ws.Range("E7").Value = "'  +0.29%  "
$ws.Range("E8").Value = "'  +0.57%  "
$ws.Range("E9").Value = "'  +0.63%  "
$ws.Range("D10").Value = "'19.17"
$ws.Range("E10").Value = "'  +1.44%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "'  +1.64%  "
$ws.Range("D12").Value = "'1.857.66"
$ws.Range("E12").Value = "'  +0.35%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.10"
$ws.Range("E13").Value = "'  +1.40%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.597.88"
$ws.Range("E14").Value = "'  -2.37%  "
$ws.Range("E15").Value = "'  -0.11%  "
$ws.Range("D16").Value = "'63.53"
$ws.Range("E16").Value = "'  +1.43%  "
$ws.Range("D17").Value = "'26.582.40"
$ws.Range("E17").Value = "'  +0.61%  "
$ws.Range("E18").Value = "'  +1.49%  "
$ws.Range("D19").Value = "'215.49"
$ws.Range("E19").Value = "'  +5.99%  "
$ws.Range("E20").Value = "'  +0.30%  "
$ws.Range("E21").Value = "'  +0.58%  "
$ws.Range("D22").Value = "'6.14"
$ws.Range("E22").Value = "'  +1.28%  "
$ws.Range("D23").Value = "'9.33"
$ws.Range("E23").Value = "'  -0.25%  "
$ws.Range("E24").Value = "'  +5.46%  "
$ws.Range("D25").Value = "'147.43"
$ws.Range("E25").Value = "'  +1.93%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "'  +0.37%  "
$ws.Range("D27").Value = "'0.120"
$ws.Range("E27").Value = "'  +0.69%  "
$ws.Range("E28").Value = "'  +3.81%  "
$ws.Range("D29").Value = "'15.49"
$ws.Range("E29").Value = "'  +1.83%  "
$ws.Range("E30").Value = "'  -2.91%  "
$ws.Range("E31").Value = "'  -0.44%  "
$ws.Range("E32").Value = "'  +2.81%  "
$ws.Range("D33").Value = "'2.96"
$ws.Range("E33").Value = "'  +1.07%  "
$ws.Range("E34").Value = "'  -0.56%  "
$ws.Range("E35").Value = "'  -0.15%  "
$ws.Range("D36").Value = "'1.220.46"
$ws.Range("E36").Value = "'  +5.42%  "
$ws.Range("E37").Value = "'  +5.13%  "
$ws.Range("D38").Value = "'0.800"
$ws.Range("E38").Value = "'  -0.71%  "
$ws.Range("E39").Value = "'  +0.29%  "
$ws.Range("E40").Value = "'  +0.11%  "
$ws.Range("D41").Value = "'2.28"
$ws.Range("E41").Value = "'  -1.77%  "
$ws.Range("D42").Value = "'0.794"
$ws.Range("E42").Value = "'  +1.52%  "
$ws.Range("D43").Value = "'5.34"
$ws.Range("E43").Value = "'  -0.92%  "
$ws.Range("D44").Value = "'1.765.53"
$ws.Range("E44").Value = "'  +0.16%  "
$ws.Range("D45").Value = "'92.99"
$ws.Range("E45").Value = "'  +1.01%  "
$ws.Range("E46").Value = "'  +2.48%  "
$ws.Range("D47").Value = "'55.05"
$ws.Range("E47").Value = "'  +1.81%  "
$ws.Range("D48").Value = "'0.0₆0103"
$ws.Range("E48").Value = "'  -0.32%  "
$ws.Range("D49").Value = "'0.0512"
$ws.Range("E49").Value = "'  +0.81%  "
$ws.Range("D50").Value = "'7.57"
$ws.Range("E50").Value = "'  +2.74%  "
$ws.Range("D51").Value = "'0.409"
$ws.Range("E51").Value = "'  -0.13%  "
